# WC.xlsx update
#  - refresh the "System Date and time" column (B2:B19) with the new run's timestamps
#  - drop the leftover blank rows (20:36) that trailed the data table
#  - reset the active selection to A2
#  - set the page to paper size 1 (Letter) on print setup

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "03-09-2024 21:57:08"
    3  = "03-09-2024 21:57:14"
    4  = "03-09-2024 21:57:20"
    5  = "03-09-2024 21:57:26"
    6  = "03-09-2024 21:57:31"
    7  = "03-09-2024 21:57:37"
    8  = "03-09-2024 21:57:43"
    9  = "03-09-2024 21:57:49"
    10 = "03-09-2024 21:57:55"
    11 = "03-09-2024 21:58:01"
    12 = "03-09-2024 21:58:06"
    13 = "03-09-2024 21:59:19"
    14 = "03-09-2024 21:59:25"
    15 = "03-09-2024 21:59:31"
    16 = "03-09-2024 21:59:37"
    17 = "03-09-2024 21:59:43"
    18 = "03-09-2024 21:59:49"
    19 = "03-09-2024 21:59:55"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("B$row").Value = $timestamps[$row]
}

# Remove the trailing empty rows (20-36), including the stray styled cell at F36.
$ws.Rows("20:36").Delete()

# Printing preference captured alongside the resave.
$ws.PageSetup.PaperSize = 1

# Cursor/selection ends up on A2 in the saved file.
$ws.Range("A2").Select() | Out-Null

Write-Output "WC.xlsx updated"
